# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Re-sort the worker mora table (rows 16-36) from grouped-by-worker (desc period)
# to grouped-by-period (asc), keeping each (worker, period) pair's Valor Mora / Salario Basico intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    @{ Row = 16; C = "73147208"; D = "GREGORIO JOSE DIAZ BANDERA"; E = "1906"; F = 64000; G = 1600000 }
    @{ Row = 17; C = "1047453567"; D = "RUBEN REYES MUÑOZ"; E = "1906"; F = 33125; G = 828116 }
    @{ Row = 18; C = "79539228"; D = "OSCAR HUMBERTO CAÑAS DIAZ"; E = "1906"; F = 64000; G = 1600000 }
    @{ Row = 19; C = "73147208"; D = "GREGORIO JOSE DIAZ BANDERA"; E = "1907"; F = 64000; G = 1600000 }
    @{ Row = 20; C = "1047453567"; D = "RUBEN REYES MUÑOZ"; E = "1907"; F = 33125; G = 828116 }
    @{ Row = 21; C = "79539228"; D = "OSCAR HUMBERTO CAÑAS DIAZ"; E = "1907"; F = 64000; G = 1600000 }
    @{ Row = 22; C = "73147208"; D = "GREGORIO JOSE DIAZ BANDERA"; E = "1908"; F = 64000; G = 1600000 }
    @{ Row = 23; C = "1047453567"; D = "RUBEN REYES MUÑOZ"; E = "1908"; F = 33125; G = 828116 }
    @{ Row = 24; C = "79539228"; D = "OSCAR HUMBERTO CAÑAS DIAZ"; E = "1908"; F = 64000; G = 1600000 }
    @{ Row = 25; C = "73147208"; D = "GREGORIO JOSE DIAZ BANDERA"; E = "1909"; F = 64000; G = 1600000 }
    @{ Row = 26; C = "1047453567"; D = "RUBEN REYES MUÑOZ"; E = "1909"; F = 33125; G = 828116 }
    @{ Row = 27; C = "79539228"; D = "OSCAR HUMBERTO CAÑAS DIAZ"; E = "1909"; F = 64000; G = 1600000 }
    @{ Row = 28; C = "73147208"; D = "GREGORIO JOSE DIAZ BANDERA"; E = "1910"; F = 64000; G = 1600000 }
    @{ Row = 29; C = "1047453567"; D = "RUBEN REYES MUÑOZ"; E = "1910"; F = 33125; G = 828116 }
    @{ Row = 30; C = "79539228"; D = "OSCAR HUMBERTO CAÑAS DIAZ"; E = "1910"; F = 64000; G = 1600000 }
    @{ Row = 31; C = "73147208"; D = "GREGORIO JOSE DIAZ BANDERA"; E = "1911"; F = 64000; G = 1600000 }
    @{ Row = 32; C = "1047453567"; D = "RUBEN REYES MUÑOZ"; E = "1911"; F = 33125; G = 828116 }
    @{ Row = 33; C = "79539228"; D = "OSCAR HUMBERTO CAÑAS DIAZ"; E = "1911"; F = 64000; G = 1600000 }
    @{ Row = 34; C = "73147208"; D = "GREGORIO JOSE DIAZ BANDERA"; E = "1912"; F = 34133; G = 1600000 }
    @{ Row = 35; C = "1047453567"; D = "RUBEN REYES MUÑOZ"; E = "1912"; F = 17667; G = 828116 }
    @{ Row = 36; C = "79539228"; D = "OSCAR HUMBERTO CAÑAS DIAZ"; E = "1912"; F = 34133; G = 1600000 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 5).Value = $item.E
    $ws.Cells.Item($item.Row, 6).Value = $item.F
    $ws.Cells.Item($item.Row, 7).Value = $item.G
}
